$d = $word.ActiveDocument

# Small helper: first paragraph whose text contains $needle.
function Find-ParagraphContaining($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- Title / cover page & body placeholders -------------------------------
# "asda" -> "DocGenerator" everywhere it is used as the product-name
# placeholder (title line + the two "... the asda tests"/"... the asda."
# sentences in the Purpose/Scope sections).
$d.Content.Find.Execute("asda", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "DocGenerator", 2)

# " asd" -> " 1.0" (version line, directly under the title)
$pVersion = Find-ParagraphContaining(" asd")
$pVersion.Range.Find.Execute(" asd", $true, $false, $false, $false, $false, `
                              $true, 1, $false, " 1.0", 2)

# --- Approval section ------------------------------------------------------
# "... approved by  asd" -> "... approved by  Eu"
$pApprover = Find-ParagraphContaining("approved by")
$pApprover.Range.Find.Execute(" asd", $true, $false, $false, $false, $false, `
                               $true, 1, $false, " Eu", 2)

# "Date:  asd" -> "Date:  11.22.2020"
$pDate = Find-ParagraphContaining("Date: ")
$pDate.Range.Find.Execute(" asd", $true, $false, $false, $false, $false, `
                           $true, 1, $false, " 11.22.2020", 2)

# --- Issues table: update the first two data rows, drop the rest ----------
$t = $d.Tables.Item(3)

$t.Cell(2, 1).Range.Text = "8391.0"
$t.Cell(2, 2).Range.Text = "Defect"
$t.Cell(2, 3).Range.Text = "Stuff happens when it should not"
$t.Cell(2, 4).Range.Text = "5.0"
$t.Cell(2, 5).Range.Text = "3.0"
$t.Cell(2, 6).Range.Text = "Open"

$t.Cell(3, 1).Range.Text = "3472.0"
$t.Cell(3, 2).Range.Text = "Story"
$t.Cell(3, 3).Range.Text = "Stuff happens when it should "
$t.Cell(3, 4).Range.Text = "7.0"
$t.Cell(3, 5).Range.Text = "6.0"
$t.Cell(3, 6).Range.Text = "Closed"

# Remove the now-unused rows (old issueType3/4/5) -- delete from the
# bottom up so earlier indices stay valid.
while ($t.Rows.Count -gt 3) {
    $t.Rows.Item($t.Rows.Count).Delete()
}
